# Atualização de bases das ligas, do dia: 15-04-2024 às 22:35
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 124 and 125: the two fixtures swap all their data (every column except
# the running index in column A, which stays put).
# ---------------------------------------------------------------------------

# New row 124 (was row 125's data)
$ws.Range("B124").Value = 7127388
$ws.Range("F124").Value = "Sydney FC"
$ws.Range("G124").Value = "Brisbane Roar"
$ws.Range("H124").Value = 1
$ws.Range("I124").Value = 1
$ws.Range("J124").Value = "D"
$ws.Range("K124").Value = 1.5
$ws.Range("L124").Value = 5
$ws.Range("M124").Value = 5
$ws.Range("N124").Value = 1.533
$ws.Range("O124").Value = 5.25
$ws.Range("P124").Value = 5
$ws.Range("Q124").Value = -1
$ws.Range("R124").Value = 1.8
$ws.Range("S124").Value = 2.05
$ws.Range("T124").Value = 3.5
$ws.Range("U124").Value = 1.925
$ws.Range("V124").Value = 1.925
$ws.Range("W124").Value = -1
$ws.Range("X124").Value = 4.25
$ws.Range("Y124").Value = -1
$ws.Range("Z124").Value = -1
$ws.Range("AA124").Value = 1.05
$ws.Range("AB124").Value = -1
$ws.Range("AC124").Value = 0.925

# New row 125 (was row 124's data)
$ws.Range("B125").Value = 7128012
$ws.Range("F125").Value = "Macarthur FC"
$ws.Range("G125").Value = "Central Coast Mariners"
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 3
$ws.Range("J125").Value = "A"
$ws.Range("K125").Value = 2.4
$ws.Range("L125").Value = 3.5
$ws.Range("M125").Value = 2.75
$ws.Range("N125").Value = 3.4
$ws.Range("O125").Value = 3.75
$ws.Range("P125").Value = 2.05
$ws.Range("Q125").Value = 0.25
$ws.Range("R125").Value = 2.025
$ws.Range("S125").Value = 1.825
$ws.Range("T125").Value = 3
$ws.Range("U125").Value = 2.05
$ws.Range("V125").Value = 1.8
$ws.Range("W125").Value = -1
$ws.Range("X125").Value = -1
$ws.Range("Y125").Value = 1.05
$ws.Range("Z125").Value = -1
$ws.Range("AA125").Value = 0.825
$ws.Range("AB125").Value = 0
$ws.Range("AC125").Value = 0

# ---------------------------------------------------------------------------
# Row 151: a fresh fixture (id 149) now carries a played match's full data
# set (final score + closing odds + P/L columns), replacing the placeholder
# fixture that used to sit here.
# ---------------------------------------------------------------------------
$ws.Range("B151").Value = 7127409
$ws.Range("E151").Value = 45396.08333333334
$ws.Range("F151").Value = "Melbourne City"
$ws.Range("G151").Value = "Perth Glory"
$ws.Range("H151").Value = 8
$ws.Range("I151").Value = 0
$ws.Range("J151").Value = "H"
$ws.Range("K151").Value = 1.571
$ws.Range("L151").Value = 4.5
$ws.Range("M151").Value = 4.75
$ws.Range("N151").Value = 1.363
$ws.Range("O151").Value = 5.25
$ws.Range("P151").Value = 7.5
$ws.Range("Q151").Value = -1.5
$ws.Range("R151").Value = 1.95
$ws.Range("S151").Value = 1.95
$ws.Range("T151").Value = 3.5
$ws.Range("U151").Value = 1.85
$ws.Range("V151").Value = 2
$ws.Range("W151").Value = 0.363
$ws.Range("X151").Value = -1
$ws.Range("Y151").Value = -1
$ws.Range("Z151").Value = 0.95
$ws.Range("AA151").Value = -1
$ws.Range("AB151").Value = 0.8500000000000001
$ws.Range("AC151").Value = -1

# ---------------------------------------------------------------------------
# Rows 152-155: each keeps the identity (match id / teams / kickoff time) of
# the fixture that used to occupy the row above it, but the odds columns get
# refreshed with newer quotes.
# ---------------------------------------------------------------------------

# Row 152 (id 150)
$ws.Range("B152").Value = 7702377
$ws.Range("E152").Value = 45398.25
$ws.Range("F152").Value = "Western United FC"
$ws.Range("G152").Value = "Adelaide United"
$ws.Range("K152").Value = 2.4
$ws.Range("L152").Value = 3.75
$ws.Range("M152").Value = 2.5
$ws.Range("N152").Value = 2.5
$ws.Range("O152").Value = 4.2
$ws.Range("P152").Value = 2.4
$ws.Range("Q152").Value = 0
$ws.Range("R152").Value = 1.99
$ws.Range("S152").Value = 1.91
$ws.Range("T152").Value = 3.5
$ws.Range("U152").Value = 1.85
$ws.Range("V152").Value = 2

# Row 153 (id 151)
$ws.Range("B153").Value = 7127410
$ws.Range("E153").Value = 45401.28125
$ws.Range("F153").Value = "Newcastle Jets"
$ws.Range("G153").Value = "Wellington Phoenix"
$ws.Range("K153").Value = 2.8
$ws.Range("L153").Value = 3.4
$ws.Range("M153").Value = 2.45
$ws.Range("N153").Value = 3
$ws.Range("O153").Value = 3.4
$ws.Range("P153").Value = 2.25
$ws.Range("Q153").Value = 0.25
$ws.Range("R153").Value = 1.88
$ws.Range("S153").Value = 2.02
$ws.Range("T153").Value = 2.75
$ws.Range("U153").Value = 1.8
$ws.Range("V153").Value = 2.05

# Row 154 (id 152)
$ws.Range("B154").Value = 8096897
$ws.Range("E154").Value = 45402.10416666666
$ws.Range("F154").Value = "Western Sydney Wanderers"
$ws.Range("G154").Value = "Melbourne City"
$ws.Range("K154").Value = 3.25
$ws.Range("L154").Value = 3.8
$ws.Range("M154").Value = 2
$ws.Range("N154").Value = 3.25
$ws.Range("O154").Value = 3.8
$ws.Range("P154").Value = 2
$ws.Range("Q154").Value = 0.5

# Row 155 (id 153)
$ws.Range("B155").Value = 7127411
$ws.Range("E155").Value = 45402.1875
$ws.Range("F155").Value = "Melbourne Victory"
$ws.Range("G155").Value = "Brisbane Roar"
$ws.Range("K155").Value = 1.65
$ws.Range("L155").Value = 4
$ws.Range("M155").Value = 4.75
$ws.Range("N155").Value = 1.65
$ws.Range("O155").Value = 4
$ws.Range("P155").Value = 4.75
$ws.Range("Q155").Value = -0.75
$ws.Range("R155").Value = 1.84
$ws.Range("S155").Value = 2.06
$ws.Range("T155").Value = 3
$ws.Range("U155").Value = 1.875
$ws.Range("V155").Value = 1.975

# ---------------------------------------------------------------------------
# Row 156 (new): the fixture that used to live in row 155 (id 154) gets
# pushed down to its own new row, also with refreshed odds. Clone the
# formatting from row 155 first (bold/bordered id cell, date-formatted
# kickoff cell) then fill in the values.
# ---------------------------------------------------------------------------
$ws.Range("A155:AC155").Copy() | Out-Null
$ws.Range("A156:AC156").PasteSpecial(-4122) | Out-Null

$ws.Range("A156").Value = 154
$ws.Range("B156").Value = 7127415
$ws.Range("C156").Value = "Australia ALeague"
$ws.Range("D156").Value = "Australia ALeague"
$ws.Range("E156").Value = 45402.28125
$ws.Range("F156").Value = "Macarthur FC"
$ws.Range("G156").Value = "Sydney FC"
$ws.Range("K156").Value = 3.8
$ws.Range("L156").Value = 4.2
$ws.Range("M156").Value = 1.8
$ws.Range("N156").Value = 4.2
$ws.Range("O156").Value = 4.2
$ws.Range("P156").Value = 1.727
$ws.Range("Q156").Value = 0.75
$ws.Range("R156").Value = 1.95
$ws.Range("S156").Value = 1.95
$ws.Range("T156").Value = 3.5
$ws.Range("U156").Value = 1.925
$ws.Range("V156").Value = 1.925
$ws.Range("W156").Value = 0
$ws.Range("X156").Value = 0
$ws.Range("Y156").Value = 0
$ws.Range("Z156").Value = 0
$ws.Range("AA156").Value = 0
